$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 666.3333
$ws.Cells.Item(28, 9).Value = 399.75
$ws.Cells.Item(28, 11).Value = 399.75
$ws.Cells.Item(28, 13).Value = 85.25

$ws.Cells.Item(33, 8).Value = 191.7
$ws.Cells.Item(33, 9).Value = 191.7
$ws.Cells.Item(33, 11).Value = 191.7
$ws.Cells.Item(33, 13).Value = 37.30000000000001

$ws.Cells.Item(53, 8).Value = 1921.2941
$ws.Cells.Item(53, 9).Value = 185
$ws.Cells.Item(53, 10).Value = 2455.5386
$ws.Cells.Item(53, 11).Value = 185
$ws.Cells.Item(53, 12).Value = 2455.5386
$ws.Cells.Item(53, 13).Value = 452
$ws.Cells.Item(53, 14).Value = -3729.5386

$ws.Cells.Item(62, 8).Value = 2777.7727
$ws.Cells.Item(62, 9).Value = 2050.9092
$ws.Cells.Item(62, 10).Value = 3504.6365
$ws.Cells.Item(62, 11).Value = 2050.9092
$ws.Cells.Item(62, 12).Value = 3504.6365
$ws.Cells.Item(62, 13).Value = -1426.9092
$ws.Cells.Item(62, 14).Value = -4752.636500000001

$ws.Cells.Item(65, 8).Value = 2777.7727
$ws.Cells.Item(65, 9).Value = 2050.9092
$ws.Cells.Item(65, 10).Value = 3504.6365
$ws.Cells.Item(65, 11).Value = 10254.546
$ws.Cells.Item(65, 12).Value = 17523.1825
$ws.Cells.Item(65, 13).Value = -7134.546
$ws.Cells.Item(65, 14).Value = -23763.1825

$ws.Cells.Item(76, 8).Value = 4276699
$ws.Cells.Item(76, 9).Value = 3218.8572
$ws.Cells.Item(76, 11).Value = 3218.8572
$ws.Cells.Item(76, 13).Value = -2903.8572

$ws.Cells.Item(79, 8).Value = 4276699
$ws.Cells.Item(79, 9).Value = 3218.8572
$ws.Cells.Item(79, 11).Value = 3218.8572
$ws.Cells.Item(79, 13).Value = -2126.8572

$ws.Cells.Item(86, 8).Value = 19991.715
$ws.Cells.Item(86, 10).Value = 19991.715
$ws.Cells.Item(86, 12).Value = 19991.715
$ws.Cells.Item(86, 14).Value = -22237.715

$ws.Cells.Item(89, 8).Value = 19991.715
$ws.Cells.Item(89, 10).Value = 19991.715
$ws.Cells.Item(89, 12).Value = 99958.575
$ws.Cells.Item(89, 14).Value = -111190.575

$ws.Cells.Item(106, 8).Value = 12347682
$ws.Cells.Item(106, 9).Value = 16668239
$ws.Cells.Item(106, 10).Value = 3234.5715
$ws.Cells.Item(106, 11).Value = 16668239
$ws.Cells.Item(106, 12).Value = 3234.5715
$ws.Cells.Item(106, 13).Value = -16667608
$ws.Cells.Item(106, 14).Value = -4496.5715

$ws.Cells.Item(116, 8).Value = 4613.3335
$ws.Cells.Item(116, 9).Value = 2365.5
$ws.Cells.Item(116, 10).Value = 6111.8887
$ws.Cells.Item(116, 11).Value = 2365.5
$ws.Cells.Item(116, 12).Value = 6111.8887
$ws.Cells.Item(116, 13).Value = 1076.5
$ws.Cells.Item(116, 14).Value = -12995.8887

$ws.Cells.Item(129, 8).Value = 271017.12
$ws.Cells.Item(129, 9).Value = 313.85715
$ws.Cells.Item(129, 11).Value = 941.5714499999999
$ws.Cells.Item(129, 13).Value = 4058.42855

$ws.Cells.Item(137, 8).Value = 1508.5186
$ws.Cells.Item(137, 9).Value = 1231.7391
$ws.Cells.Item(137, 11).Value = 3695.2173
$ws.Cells.Item(137, 13).Value = -1145.2173

$ws.Cells.Item(138, 8).Value = 1913.0521
$ws.Cells.Item(138, 9).Value = 1377.3214
$ws.Cells.Item(138, 10).Value = 2133.647
$ws.Cells.Item(138, 11).Value = 4131.9642
$ws.Cells.Item(138, 12).Value = 6400.941
$ws.Cells.Item(138, 13).Value = 1008.0358
$ws.Cells.Item(138, 14).Value = -16680.941

$ws.Cells.Item(141, 8).Value = 1791.6522
$ws.Cells.Item(141, 9).Value = 1558.5
$ws.Cells.Item(141, 11).Value = 4675.5
$ws.Cells.Item(141, 13).Value = 504.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3772.6904
$ws.Cells.Item(32, 9).Value = 3124.077
$ws.Cells.Item(32, 10).Value = 12204.667
$ws.Cells.Item(32, 11).Value = 3124.077
$ws.Cells.Item(32, 12).Value = 12204.667
$ws.Cells.Item(32, 13).Value = -2837.077
$ws.Cells.Item(32, 14).Value = -12778.667

$ws.Cells.Item(63, 8).Value = 15626142
$ws.Cells.Item(63, 9).Value = 2285
$ws.Cells.Item(63, 11).Value = 2285
$ws.Cells.Item(63, 13).Value = -1599

$ws.Cells.Item(66, 8).Value = 15626142
$ws.Cells.Item(66, 9).Value = 2285
$ws.Cells.Item(66, 11).Value = 11425
$ws.Cells.Item(66, 13).Value = -7993

$ws.Cells.Item(74, 8).Value = 125001630
$ws.Cells.Item(74, 9).Value = 200000910
$ws.Cells.Item(74, 11).Value = 200000910
$ws.Cells.Item(74, 13).Value = -200000036

$ws.Cells.Item(77, 8).Value = 125001630
$ws.Cells.Item(77, 9).Value = 200000910
$ws.Cells.Item(77, 11).Value = 1000004550
$ws.Cells.Item(77, 13).Value = -1000000182

$ws.Cells.Item(94, 8).Value = 31165
$ws.Cells.Item(94, 9).Value = 22000
$ws.Cells.Item(94, 10).Value = 40330
$ws.Cells.Item(94, 11).Value = 22000
$ws.Cells.Item(94, 12).Value = 40330
$ws.Cells.Item(94, 13).Value = -21099
$ws.Cells.Item(94, 14).Value = -42132

$ws.Cells.Item(102, 8).Value = 2500
$ws.Cells.Item(102, 9).Value = 2500
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 2500
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -878
$ws.Cells.Item(102, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 12192.341
$ws.Cells.Item(132, 9).Value = 1386.0555
$ws.Cells.Item(132, 10).Value = 47558.363
$ws.Cells.Item(132, 11).Value = 4158.166499999999
$ws.Cells.Item(132, 12).Value = 142675.089
$ws.Cells.Item(132, 13).Value = -1628.166499999999
$ws.Cells.Item(132, 14).Value = -147735.089

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2542.8572
$ws.Cells.Item(99, 9).Value = 1950
$ws.Cells.Item(99, 10).Value = 3333.3333
$ws.Cells.Item(99, 11).Value = 1950
$ws.Cells.Item(99, 12).Value = 3333.3333
$ws.Cells.Item(99, 13).Value = -452
$ws.Cells.Item(99, 14).Value = -6329.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 13740.207
$ws.Cells.Item(31, 9).Value = 23142.072
$ws.Cells.Item(31, 10).Value = 4965.1333
$ws.Cells.Item(31, 11).Value = 23142.072
$ws.Cells.Item(31, 12).Value = 4965.1333
$ws.Cells.Item(31, 13).Value = -22847.072
$ws.Cells.Item(31, 14).Value = -5555.1333

$ws.Cells.Item(34, 8).Value = 13740.207
$ws.Cells.Item(34, 9).Value = 23142.072
$ws.Cells.Item(34, 10).Value = 4965.1333
$ws.Cells.Item(34, 11).Value = 23142.072
$ws.Cells.Item(34, 12).Value = 4965.1333
$ws.Cells.Item(34, 13).Value = -22940.072
$ws.Cells.Item(34, 14).Value = -5369.1333

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(57, 8).Value = 6333.3335
$ws.Cells.Item(57, 9).Value = 1000
$ws.Cells.Item(57, 10).Value = 9000
$ws.Cells.Item(57, 11).Value = 3000
$ws.Cells.Item(57, 12).Value = 27000
$ws.Cells.Item(57, 13).Value = -2441
$ws.Cells.Item(57, 14).Value = -28118

$ws.Cells.Item(113, 8).Value = 534.17645
$ws.Cells.Item(113, 9).Value = 453.77777
$ws.Cells.Item(113, 10).Value = 624.625
$ws.Cells.Item(113, 11).Value = 1361.33331
$ws.Cells.Item(113, 12).Value = 1873.875
$ws.Cells.Item(113, 13).Value = 808.66669
$ws.Cells.Item(113, 14).Value = -6213.875

$ws.Cells.Item(129, 8).Value = 626450
$ws.Cells.Item(129, 9).Value = 800
$ws.Cells.Item(129, 10).Value = 835000
$ws.Cells.Item(129, 11).Value = 2400
$ws.Cells.Item(129, 12).Value = 2505000
$ws.Cells.Item(129, 13).Value = 2600
$ws.Cells.Item(129, 14).Value = -2515000

$ws.Cells.Item(131, 8).Value = 743.9299999999999
$ws.Cells.Item(131, 10).Value = 743.9299999999999
$ws.Cells.Item(131, 12).Value = 2231.79
$ws.Cells.Item(131, 14).Value = -12311.79

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 25000
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 14).ClearContents()

$ws.Cells.Item(57, 8).Value = 27930
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 10).Value = 27930
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 12).Value = 27930
$ws.Cells.Item(57, 13).ClearContents()
$ws.Cells.Item(57, 14).Value = -29570

$ws.Cells.Item(80, 8).Value = 3482
$ws.Cells.Item(80, 9).Value = 3245
$ws.Cells.Item(80, 10).Value = 3640
$ws.Cells.Item(80, 11).Value = 3245
$ws.Cells.Item(80, 12).Value = 3640
$ws.Cells.Item(80, 13).Value = -2247
$ws.Cells.Item(80, 14).Value = -5636

$ws.Cells.Item(83, 8).Value = 3482
$ws.Cells.Item(83, 9).Value = 3245
$ws.Cells.Item(83, 10).Value = 3640
$ws.Cells.Item(83, 11).Value = 16225
$ws.Cells.Item(83, 12).Value = 18200
$ws.Cells.Item(83, 13).Value = -11233
$ws.Cells.Item(83, 14).Value = -28184

$ws.Cells.Item(132, 8).Value = 33123.39
$ws.Cells.Item(132, 9).Value = 5623.0713
$ws.Cells.Item(132, 10).Value = 129374.5
$ws.Cells.Item(132, 11).Value = 16869.2139
$ws.Cells.Item(132, 12).Value = 388123.5
$ws.Cells.Item(132, 13).Value = -14339.2139
$ws.Cells.Item(132, 14).Value = -393183.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 3069.4
$ws.Cells.Item(93, 9).Value = 2711.25
$ws.Cells.Item(93, 10).Value = 4502
$ws.Cells.Item(93, 11).Value = 2711.25
$ws.Cells.Item(93, 12).Value = 4502
$ws.Cells.Item(93, 13).Value = -1463.25
$ws.Cells.Item(93, 14).Value = -6998

$ws.Cells.Item(100, 8).Value = 2952.9412
$ws.Cells.Item(100, 9).Value = 2624.75
$ws.Cells.Item(100, 10).Value = 3053.923
$ws.Cells.Item(100, 11).Value = 2624.75
$ws.Cells.Item(100, 12).Value = 3053.923
$ws.Cells.Item(100, 13).Value = -2083.75
$ws.Cells.Item(100, 14).Value = -4135.923

$ws.Cells.Item(132, 8).Value = 1952.2106
$ws.Cells.Item(132, 9).Value = 1246.2667
$ws.Cells.Item(132, 10).Value = 4599.5
$ws.Cells.Item(132, 11).Value = 3738.800099999999
$ws.Cells.Item(132, 12).Value = 13798.5
$ws.Cells.Item(132, 13).Value = -1208.800099999999
$ws.Cells.Item(132, 14).Value = -18858.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4735
$ws.Cells.Item(62, 9).Value = 3504.2
$ws.Cells.Item(62, 10).Value = 5294.4546
$ws.Cells.Item(62, 11).Value = 3504.2
$ws.Cells.Item(62, 12).Value = 5294.4546
$ws.Cells.Item(62, 13).Value = -2880.2
$ws.Cells.Item(62, 14).Value = -6542.4546

$ws.Cells.Item(65, 8).Value = 4735
$ws.Cells.Item(65, 9).Value = 3504.2
$ws.Cells.Item(65, 10).Value = 5294.4546
$ws.Cells.Item(65, 11).Value = 17521
$ws.Cells.Item(65, 12).Value = 26472.273
$ws.Cells.Item(65, 13).Value = -14401
$ws.Cells.Item(65, 14).Value = -32712.273

$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(129, 14).ClearContents()
